$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(22.7957292726789, 19.22985656262092, 26.7291141594762),
    @(22.01324122430034, 4.985943378995435, 40.23106118721454),
    @(22.17683594497943, 12.93579588280059, 31.8128034117784),
    @(26.62176316608863, 8.391024657534235, 45.77450913242003),
    @(30.22840946540404, 23.85192929895855, 35.90123295683938),
    @(20.34688352132542, 17.43954552090043, 23.16562185394933),
    @(3.747721289597126, 2.893905736009925, 5.116101568629986),
    @(28.08733651613845, 23.44276965691799, 33.06353489111826),
    @(3.767801038626333, 2.89233431784446, 5.14548468662595),
    @(15.41814243619785, 11.23587105155694, 20.4678900674466),
    @(32.38415097798217, 28.81439764502737, 36.31146989734884),
    @(5.65157898956937, 4.120185087888845, 7.925955876818584),
    @(9.122950656112311, 6.595662102864509, 12.78476494905586),
    @(32.23696667293063, 10.47129575342465, 54.48023090345728),
    @(30.34993155775826, 23.92223098370977, 36.02864779402602),
    @(32.20815545889828, 28.76636088901758, 36.08579470249187),
    @(22.17683594497943, 12.93579588280059, 31.8128034117784),
    @(15.75624207998909, 11.50318374142364, 21.9119032254371),
    @(35.27651753834856, 29.67331660380743, 40.36306633673983),
    @(32.68772957053474, 21.1448173133303, 43.61960605546841),
    @(3.824864494852398, 2.913011296233764, 5.233377624053453),
    @(34.56891033246549, 9.172672716894972, 57.14304279191119),
    @(30.33912896354122, 23.96240590232645, 35.98415916405936),
    @(25.41795868387253, 20.07533898758255, 29.80032808673684),
    @(32.96606823938602, 22.04653618721458, 43.9567899630034),
    @(5.59677433792848, 4.109584185395674, 8.094068839418723),
    @(21.68448429728166, 18.76509775256399, 24.8218263675423),
    @(36.34420247164773, 11.12510570776251, 57.31924722765815),
    @(29.97245623510285, 26.93860789202504, 33.32993014576635),
    @(9.340605318891683, 6.56517562034421, 13.1169090752892),
    @(5.524349790936562, 4.100575966059266, 7.867658177867092),
    @(9.269542438902693, 6.706173066015029, 12.95023602790925),
    @(25.74718358503872, 20.96019774114424, 31.38176446472727),
    @(32.51572561956845, 26.74853824061113, 37.73993448926148),
    @(20.36972614876687, 17.40272070592649, 23.20212218093124),
    @(36.19116775332004, 23.57124184959178, 48.45225262672845),
    @(9.177237881435854, 6.556961302997985, 13.296856322919),
    @(33.5356449767108, 29.52652888212343, 37.62347740533664),
    @(26.680144110805, 13.97311799266535, 38.04235071902983),
    @(18.93162566153917, 16.73758577934632, 21.46885416020572),
    @(15.62966812748561, 11.50295470583448, 20.78064197972285),
    @(29.97365148303983, 26.93589379001958, 33.33764214153383),
    @(13.83128420345437, 10.89918134499355, 17.1341447355329),
    @(15.96709044118684, 11.75487283392216, 21.51739349210451),
    @(8.90750635555338, 6.462626924795975, 12.21514508074283),
    @(3.795284418311202, 2.834920445855212, 5.276497808235737),
    @(32.53326962256635, 28.93779961626342, 36.47742860563646),
    @(25.2478975674148, 15.27996804273653, 35.98422863992165),
    @(36.17365981974445, 23.58432854684905, 48.37888460195174),
    @(21.97260688674001, 5.303065296803655, 40.16996529680357),
    @(8.820770231815736, 6.532839144147228, 12.20975019517104),
    @(33.64028560162796, 29.75870271071227, 37.73851917428607),
    @(32.06032767406005, 9.472675547945206, 54.51456710697978),
    @(32.79752668842855, 20.77501086385707, 44.3790375009652),
    @(9.161590642464779, 6.668565714863298, 12.88963408781841),
    @(5.550337331710252, 4.095500327763819, 7.879327174488426),
    @(36.07043797559354, 23.81036574864756, 49.07575245348662),
    @(14.9844822381122, 11.29970321402797, 19.14078396781226),
    @(35.20907594237234, 11.56186906392693, 56.7927564253098),
    @(9.162047595752968, 6.629833726366762, 12.87679342816766),
    @(15.43582356709032, 11.29250027079209, 20.45506690956195),
    @(9.19330025705392, 6.465816138203738, 12.57386067816082),
    @(35.21484550884096, 11.0447670091324, 56.60139947814737),
    @(5.525844263879763, 4.079875060133376, 7.652851166556045),
    @(24.88070292301152, 20.22248432702806, 30.10838592211701),
    @(23.96375686257021, 20.01131704059455, 27.98742621818564),
    @(30.02401447478087, 18.2541578071959, 40.93696674014414),
    @(3.827299385830869, 2.8874114370634, 5.179403077591576),
    @(26.77615651454181, 15.41291650385067, 37.04061809375549),
    @(21.8783063315752, 19.17187175290154, 24.60168772840371),
    @(26.35490243449845, 23.43495790335223, 29.32852155874417)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $values[$i][0]
    $ws.Cells.Item($r, 2).Value = $values[$i][1]
    $ws.Cells.Item($r, 3).Value = $values[$i][2]
}